$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2043.1875
$ws.Range("J17").Value = 2168.5386
$ws.Range("L17").Value = 6505.6158
$ws.Range("N17").Value = -6841.6158

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 50272.055
$ws.Range("I28").Value = 242.8125
$ws.Range("K28").Value = 242.8125
$ws.Range("M28").Value = 242.1875

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 12566.667
$ws.Range("I40").Value = 6042.857
$ws.Range("J40").Value = 16718.182
$ws.Range("K40").Value = 6042.857
$ws.Range("L40").Value = 16718.182
$ws.Range("M40").Value = -5867.857
$ws.Range("N40").Value = -17068.182

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1104.6666
$ws.Range("I58").Value = 85.7
$ws.Range("J58").Value = 6199.5
$ws.Range("K58").Value = 257.1
$ws.Range("L58").Value = 18598.5
$ws.Range("M58").Value = -107.1
$ws.Range("N58").Value = -18898.5

# ALC row 61
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 220.14285
$ws.Range("I61").Value = 220.14285
$ws.Range("K61").Value = 660.4285500000001
$ws.Range("M61").Value = -488.4285500000001

# ALC row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 9617.333000000001
$ws.Range("I94").Value = 9617.333000000001
$ws.Range("K94").Value = 9617.333000000001
$ws.Range("M94").Value = -9166.333000000001

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 681.9
$ws.Range("I96").Value = 427.16666
$ws.Range("J96").Value = 1064
$ws.Range("K96").Value = 1281.49998
$ws.Range("L96").Value = 3192
$ws.Range("M96").Value = 91.50001999999995
$ws.Range("N96").Value = -5938

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3066
$ws.Range("J100").Value = 3999
$ws.Range("L100").Value = 3999
$ws.Range("N100").Value = -5081

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 462999.9
$ws.Range("I137").Value = 1659.6
$ws.Range("J137").Value = 895506.4399999999
$ws.Range("K137").Value = 4978.799999999999
$ws.Range("L137").Value = 2686519.32
$ws.Range("M137").Value = -2428.799999999999
$ws.Range("N137").Value = -2691619.32

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 12659977
$ws.Range("I138").Value = 1011.7917
$ws.Range("J138").Value = 18183890
$ws.Range("K138").Value = 3035.3751
$ws.Range("L138").Value = 54551670
$ws.Range("M138").Value = 2104.6249
$ws.Range("N138").Value = -54561950

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9301.852999999999
$ws.Range("I32").Value = 6052.2183
$ws.Range("K32").Value = 6052.2183
$ws.Range("M32").Value = -5765.2183

# ARM row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 55333.332
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 55333.332
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 55333.332
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -55963.332

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1832.85
$ws.Range("I61").Value = 1421
$ws.Range("J61").Value = 4166.6665
$ws.Range("K61").Value = 1421
$ws.Range("L61").Value = 4166.6665
$ws.Range("M61").Value = -1209
$ws.Range("N61").Value = -4590.6665

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 38701.125
$ws.Range("I102").Value = 18264.666
$ws.Range("J102").Value = 100010.5
$ws.Range("K102").Value = 18264.666
$ws.Range("L102").Value = 100010.5
$ws.Range("M102").Value = -16642.666
$ws.Range("N102").Value = -103254.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1832.85
$ws.Range("I136").Value = 1421
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 4263
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -1713
$ws.Range("N136").Value = -17599.9995

# ARM row 138
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 113142.664
$ws.Range("J138").Value = 113142.664
$ws.Range("L138").Value = 113142.664
$ws.Range("N138").Value = -123422.664

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 85004.836
$ws.Range("J139").Value = 85004.836
$ws.Range("L139").Value = 85004.836
$ws.Range("N139").Value = -95284.836

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 112435.15
$ws.Range("I20").Value = 148711.67
$ws.Range("J20").Value = 3605.6
$ws.Range("K20").Value = 148711.67
$ws.Range("L20").Value = 3605.6
$ws.Range("M20").Value = -148464.67
$ws.Range("N20").Value = -4099.6

# BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 34396.715
$ws.Range("I81").Value = 20000
$ws.Range("J81").Value = 36796.168
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 36796.168
$ws.Range("M81").Value = -18939
$ws.Range("N81").Value = -38918.168

# BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 34396.715
$ws.Range("I84").Value = 20000
$ws.Range("J84").Value = 36796.168
$ws.Range("K84").Value = 60000
$ws.Range("L84").Value = 110388.504
$ws.Range("M84").Value = -54696
$ws.Range("N84").Value = -120996.504

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 885.62067
$ws.Range("I94").Value = 826.2941
$ws.Range("J94").Value = 969.6667
$ws.Range("K94").Value = 826.2941
$ws.Range("L94").Value = 969.6667
$ws.Range("M94").Value = -375.2941
$ws.Range("N94").Value = -1871.6667

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1490186.6
$ws.Range("I99").Value = 1348.6364
$ws.Range("J99").Value = 3127908.5
$ws.Range("K99").Value = 1348.6364
$ws.Range("L99").Value = 3127908.5
$ws.Range("M99").Value = 149.3635999999999
$ws.Range("N99").Value = -3130904.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3818.0789
$ws.Range("J31").Value = 5376.5293
$ws.Range("L31").Value = 5376.5293
$ws.Range("N31").Value = -5966.5293

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3818.0789
$ws.Range("J34").Value = 5376.5293
$ws.Range("L34").Value = 5376.5293
$ws.Range("N34").Value = -5780.5293

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1843.0741
$ws.Range("I132").Value = 926
$ws.Range("K132").Value = 2778
$ws.Range("M132").Value = -248

# CRP row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 62000
$ws.Range("J135").Value = 62000
$ws.Range("L135").Value = 62000
$ws.Range("N135").Value = -72140

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 14.1
$ws.Range("I12").Value = 15.75
$ws.Range("J12").Value = 13
$ws.Range("K12").Value = 47.25
$ws.Range("L12").Value = 39
$ws.Range("M12").Value = 125.75
$ws.Range("N12").Value = -385

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1072.375
$ws.Range("I68").Value = 1195.5
$ws.Range("J68").Value = 949.25
$ws.Range("K68").Value = 3586.5
$ws.Range("L68").Value = 2847.75
$ws.Range("M68").Value = -2775.5
$ws.Range("N68").Value = -4469.75

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1072.375
$ws.Range("I71").Value = 1195.5
$ws.Range("J71").Value = 949.25
$ws.Range("K71").Value = 10759.5
$ws.Range("L71").Value = 8543.25
$ws.Range("M71").Value = -6703.5
$ws.Range("N71").Value = -16655.25

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 884338.1
$ws.Range("I121").Value = 364.57144
$ws.Range("K121").Value = 1093.71432
$ws.Range("M121").Value = 216.28568

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 634444.25
$ws.Range("I70").Value = 615714.1
$ws.Range("K70").Value = 615714.1
$ws.Range("M70").Value = -615444.1

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 634444.25
$ws.Range("I73").Value = 615714.1
$ws.Range("K73").Value = 615714.1
$ws.Range("M73").Value = -614778.1

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2248.75
$ws.Range("I132").Value = 1808.3334
$ws.Range("K132").Value = 5425.0002
$ws.Range("M132").Value = -2895.0002

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5764520.5
$ws.Range("I40").Value = 2904.3333
$ws.Range("K40").Value = 2904.3333
$ws.Range("M40").Value = -2768.3333

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 3127167.5
$ws.Range("I55").Value = 1198.4546
$ws.Range("J55").Value = 10004299
$ws.Range("K55").Value = 1198.4546
$ws.Range("L55").Value = 10004299
$ws.Range("M55").Value = -1025.4546
$ws.Range("N55").Value = -10004645

# LTW row 88
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 32999.168
$ws.Range("I88").Value = 24499.5
$ws.Range("K88").Value = 24499.5
$ws.Range("M88").Value = -24071.5

# LTW row 91
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H91").Value = 32999.168
$ws.Range("I91").Value = 24499.5
$ws.Range("K91").Value = 24499.5
$ws.Range("M91").Value = -23017.5

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1902.5625
$ws.Range("I93").Value = 1538.6522
$ws.Range("K93").Value = 1538.6522
$ws.Range("M93").Value = -290.6522

# LTW row 96
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 109999.5
$ws.Range("J96").Value = 109999.5
$ws.Range("L96").Value = 109999.5
$ws.Range("N96").Value = -115491.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9392210
$ws.Range("I122").Value = 22459.521
$ws.Range("K122").Value = 67378.56299999999
$ws.Range("M122").Value = -64928.56299999999

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 79988
$ws.Range("J133").Value = 79988
$ws.Range("L133").Value = 79988
$ws.Range("N133").Value = -85048

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6104.758
$ws.Range("I136").Value = 7052.2104
$ws.Range("K136").Value = 21156.6312
$ws.Range("M136").Value = -18606.6312

# LTW row 141
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 112814
$ws.Range("J141").Value = 112814
$ws.Range("L141").Value = 112814
$ws.Range("N141").Value = -123174

# WVR row 44
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 41499.5
$ws.Range("J44").Value = 41499.5
$ws.Range("L44").Value = 41499.5
$ws.Range("N44").Value = -42607.5

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 479.55
$ws.Range("I113").Value = 295.7143
$ws.Range("J113").Value = 908.5
$ws.Range("K113").Value = 887.1428999999999
$ws.Range("L113").Value = 2725.5
$ws.Range("M113").Value = 1282.8571
$ws.Range("N113").Value = -7065.5

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2944.487
$ws.Range("I122").Value = 1606.0385
$ws.Range("J122").Value = 5621.385
$ws.Range("K122").Value = 4818.1155
$ws.Range("L122").Value = 16864.155
$ws.Range("M122").Value = -2368.1155
$ws.Range("N122").Value = -21764.155

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2399.9807
$ws.Range("I136").Value = 2076.5
$ws.Range("K136").Value = 6229.5
$ws.Range("M136").Value = -3679.5
